$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Paragraph: "There's a lot of people..." -> "There's an important difference..." ---
Replace-Text `
    "There’s a lot of people that think that pointing out the differences between" `
    "There’s an important difference between the hard sciences and the"

Replace-Text `
    " pseudo-sciences from the sciences is somehow denigrating the respectability of the pseudo-sciences. I disagree. " `
    " pseudo-sciences called falsifiability. "

# --- Paragraph: "The term pseudo-science was coined..." -> add "Austrian" ---
Replace-Text `
    " was coined by a 20th century philosopher, Karl Popper. " `
    " was coined by a 20th century Austrian philosopher, Karl Popper. "

# --- Paragraph: "Popper studied people like Einstein and Freud..." -> "supported" ---
Replace-Text `
    "looked for evidence that would disprove their theories, but the Freuds only looked for evidence that would support their theories." `
    "looked for evidence that would disprove their theories, but the Freuds only looked for evidence that supported their theories."

# --- Paragraph: "A scientist's best hypotheses..." ---
Replace-Text `
    "A scientist’s best hypotheses and theories are always tentative because some unthought-of experiment or a new piece of evidence could always prove them false. " `
    "A scientist’s best hypotheses and theories are always tentative because it just takes one person with an unthought-of experiment or a new piece of evidence to falsify a claim. "

Replace-Text `
    " pseudo-scientist’s theories is true as soon as " `
    " pseudo-scientist’s theories are true as soon as "

# --- Paragraph: "In other words, Einstein's Theory of Relativity..." ---
Replace-Text `
    "person present some evidence" `
    "person presents some evidence"

Replace-Text `
    "no one tested before. How do you " `
    "no one’s tested before. How do you "

# --- Paragraph: "Operationalizing lets scientists weed themselves out..." ---
Replace-Text `
    "Operationalizing lets scientists weed themselves out from the posers. Claims that’ve survived persistent attempts to disprove them automatically prompt credibility. Popper (a pseudo-scientist himself) distinguished pseudo-science from science to find better logical justifications for knowledge claims—not to dismiss " `
    "Operationalizing lets scientists weed themselves out from the posers. Claims that can survive persistent attempts to disprove them automatically garner credibility. Popper (a pseudo-scientist himself) distinguished pseudo-science from science to find better logical justifications for knowledge claims—not to dismiss "

Write-Output "done"
